# Quarterly indexing esoteric bug-fix operation
#
# Column A (rows 2-73) holds Excel date-serial values that were incorrectly
# set to the 1st of a quarter-start month. The fix re-indexes each date to
# the 15th of the following month (i.e. the mid-point of the 2nd month of
# the quarter), which is the intended "quarter of quarter" anchor date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel's (1900 date system) day-zero reference point.
$epoch = Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0

for ($r = 2; $r -le 73; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $oldVal = [double]$cell.Value2

    # Turn the existing serial date into a real date.
    $d = $epoch.AddDays($oldVal)

    # Move forward one month, keeping the year rollover correct.
    $y = $d.Year
    $m = $d.Month + 1
    if ($m -gt 12) {
        $m = $m - 12
        $y = $y + 1
    }

    # Re-anchor on the 15th of that following month.
    $next = Get-Date -Year $y -Month $m -Day 15 -Hour 0 -Minute 0 -Second 0

    # Write back as an Excel serial date (OLE Automation date == Excel's 1900 system).
    $cell.Value2 = $next.ToOADate()
}
